$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'319.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'3.43%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'1.03%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'2.36%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07717"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.17%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.685"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.88%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9441"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.75%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Value = "'-0.95%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.01%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09202"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.34%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04120"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.02%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.25%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001287"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.04%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.006042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'6.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E17").Value = "'0.04%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'4.342"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.27%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3353"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.06%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.404"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'21.11%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-2.75%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-0.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04037"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001272"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.01%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004120"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.18%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02538"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'4.79%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05324"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.68%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007785"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.50%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'1.01%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007046"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'3.46%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002154"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'15.66%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008304"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'11.80%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3482"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.12%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006682"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.72%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.09%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.2003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'55.64%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'40.20%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.09%"
$ws.Range("E51").Style = "Normal"
